# FFL_Data.xlsx - Week 7 results + team name changes
#
# 1. Append the Week 7 results (rows 74-85) including the Image column.
# 2. Clear the Week 6 "Image" column (the Image column only ever holds the
#    most-recently-completed week's team icons).
# 3. Rename teams that carried over from Week 1 through Week 7:
#      "You Gotta Gibbs!"    -> "Golden Knights"
#      "Show Me the Mooney"  -> "Moonies of Io"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Add Week 7 results (still using the old team names) ---------------
$week7 = @(
    @("Quad Goals",                     "Week 7", 68.959999999999994, 115.64,             "Team Icons/quad-modified.png"),
    @("Egbuka di Beppo",                "Week 7", 186.52,             88.24,               "Team Icons/egbuka-modified.png"),
    @("CoHo Chicken Tet",               "Week 7", 88.24,              186.52,              "Team Icons/coho-modified.png"),
    @("Stone (injured) Kittles",        "Week 7", 83.3,               67.22,               "Team Icons/stone-modified.png"),
    @("The Legend of Drewkeys",         "Week 7", 71.540000000000006, 81.180000000000007,  "Team Icons/drewkeys-modified.png"),
    @("You Gotta Gibbs!",               "Week 7", 166.86,             114.18,              "Team Icons/gibbs-modified.png"),
    @("Show Me the Mooney",             "Week 7", 131.74,             80.92,               "Team Icons/mooney-modified.png"),
    @("Bucky Charms",                   "Week 7", 115.64,             68.959999999999994,  "Team Icons/charms-modified.png"),
    @("he was a skattebo",              "Week 7", 114.18,             166.86,              "Team Icons/skattebo-modified.png"),
    @("McConkey Kong Country",          "Week 7", 80.92,              131.74,              "Team Icons/mcconkey-modified.png"),
    @("UNLIMITED BOWERS",               "Week 7", 67.22,              83.3,                "Team Icons/unlimited-modified.png"),
    @("Fantasy Champ 2022 and 2024",    "Week 7", 81.180000000000007, 71.540000000000006,  "Team Icons/ffc-modified.png")
)

$startRow = 74
for ($i = 0; $i -lt $week7.Length; $i++) {
    $r = $startRow + $i
    $data = $week7[$i]
    $ws.Cells.Item($r, 1).Value2 = $data[0]
    $ws.Cells.Item($r, 2).Value2 = $data[1]
    $ws.Cells.Item($r, 3).Value2 = $data[2]
    $ws.Cells.Item($r, 4).Value2 = $data[3]
    $ws.Cells.Item($r, 5).Value2 = $data[4]
}

# --- 2. Clear last week's (Week 6) Image column ----------------------------
$ws.Range("E62:E73").ClearContents() | Out-Null

# --- 3. Team renames (affects every week row + any other references) ------
$ws.Columns.Item(1).Replace("You Gotta Gibbs!", "Golden Knights") | Out-Null
$ws.Columns.Item(1).Replace("Show Me the Mooney", "Moonies of Io") | Out-Null

# --- 4. Drop a few stray empty-but-styled cells below the table -----------
$ws.Range("C98:C100").Clear() | Out-Null

# --- 5. Leave the selection where the author left it on save --------------
$ws.Range("E12").Select() | Out-Null
